# Actualización automática 2025-06-02 13:21:56
# Adds a new "PRESUPUESTO" (budget) column G to the "VENTA MENSUAL" sheet,
# filled with 0 values for every data row and the running total row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# --- Header cell G1 -------------------------------------------------------
$ws.Range("G1").Value = "PRESUPUESTO"
$ws.Range("G1").Font.Bold = $true
$ws.Range("G1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("G1").VerticalAlignment = -4160     # xlTop
$ws.Range("G1").Borders.LineStyle = 1
$ws.Range("G1").Borders.Weight = 2

# --- Data rows G2:G29 -------------------------------------------------------
$ws.Range("G2:G29").Value = 0
$ws.Range("G2:G29").NumberFormat = """$""#,##0.00"

# --- Total row G30 ---------------------------------------------------------
$ws.Range("G30").Value = 0
$ws.Range("G30").NumberFormat = """$""#,##0.00"
$ws.Range("G30").HorizontalAlignment = -4152  # xlRight

# --- Column width ------------------------------------------------------
$ws.Columns.Item(7).ColumnWidth = 16.17
